$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44214
$ws.Range("J2").Value = 7000

# Row 3
$ws.Range("D3").Value = 44230
$ws.Range("J3").Value = 16000
$ws.Range("K3").Value = 3000
$ws.Range("M3").Value = 3000
$ws.Range("P3").Value = 30

# Row 4
$ws.Range("D4").Value = 44600
$ws.Range("J4").Value = 1300
$ws.Range("K4").Value = 3500
$ws.Range("L4").Value = 4000
$ws.Range("M4").Value = 3808
$ws.Range("O4").Value = 'Región Metropolitana'
$ws.Range("P4").Value = 38

# Row 5
$ws.Range("D5").Value = 44215
$ws.Range("J5").Value = 16000

# Row 6
$ws.Range("D6").Value = 44229
$ws.Range("J6").Value = 16000

# Row 7
$ws.Range("D7").Value = 44210
$ws.Range("J7").Value = 8800
$ws.Range("K7").Value = 2500
$ws.Range("M7").Value = 2750
$ws.Range("P7").Value = 28

# Row 8
$ws.Range("D8").Value = 44159
$ws.Range("J8").Value = 7000

# Row 9
$ws.Range("D9").Value = 44181
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 12000
$ws.Range("K9").Value = 3000
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 3000
$ws.Range("P9").Value = 30

# Row 10
$ws.Range("D10").Value = 44186
$ws.Range("J10").Value = 10000

# Row 11
$ws.Range("D11").Value = 44167
$ws.Range("J11").Value = 7000

# Row 12
$ws.Range("D12").Value = 44188
$ws.Range("J12").Value = 12000

# Row 13
$ws.Range("D13").Value = 44166
$ws.Range("J13").Value = 7000
$ws.Range("K13").Value = 3000
$ws.Range("M13").Value = 3000
$ws.Range("P13").Value = 30

# Row 14
$ws.Range("D14").Value = 44162

# Row 15
$ws.Range("D15").Value = 44232

# Row 16
$ws.Range("D16").Value = 44161
$ws.Range("J16").Value = 7000

# Row 17
$ws.Range("D17").Value = 44189
$ws.Range("J17").Value = 16000

# Row 18
$ws.Range("D18").Value = 44160
$ws.Range("J18").Value = 7000
$ws.Range("K18").Value = 3000
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = 3000
$ws.Range("O18").Value = 'Provincia de Chacabuco'
$ws.Range("P18").Value = 30

# Row 19
$ws.Range("D19").Value = 44204
$ws.Range("J19").Value = 7000

# Row 20
$ws.Range("D20").Value = 44209
$ws.Range("K20").Value = 2500
$ws.Range("M20").Value = 2750
$ws.Range("P20").Value = 28

# Row 21
$ws.Range("D21").Value = 44231
$ws.Range("J21").Value = 12000

# Row 22
$ws.Range("D22").Value = 44168
$ws.Range("J22").Value = 7000

# Row 23
$ws.Range("D23").Value = 44187
$ws.Range("J23").Value = 12000

# Row 24
$ws.Range("D24").Value = 44602
$ws.Range("J24").Value = 12000
$ws.Range("O24").Value = 'Provincia de Chacabuco'

# Row 25
$ws.Range("D25").Value = 44602
$ws.Range("J25").Value = 6000
$ws.Range("O25").Value = 'Provincia de Chacabuco'

# Row 26
$ws.Range("D26").Value = 44245
$ws.Range("J26").Value = 9000
$ws.Range("O26").Value = 'Región Metropolitana'

# Row 27
$ws.Range("D27").Value = 44245
$ws.Range("I27").Value = 'Segunda'
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 2500
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = 2500
$ws.Range("O27").Value = 'Región Metropolitana'
$ws.Range("P27").Value = 25
